# "Update contact log template" -
# - Retitle the merged header cell from the client-specific
#   "Contact Log : Christopher Martens Law Corp" to the generic "Contact Log".
# - Leave the cursor/selection on the header row (A1:D1) instead of D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Contact Log"
$ws.Range("A1:D1").Select()
